# Activate wetland 10 (IDs 33827, 33807; ~21 acres adjacent to Springfield):
# append the new "Demo_Baseline_2010-18 C600" run result to the results table
# on the "2010 and 2010-18" sheet, right after the existing "CW3M 1.1.0 /
# Demo_Baseline_2010-19 C564" row, pushing the "2010-19" summary block down.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2010 and 2010-18")

# Insert 5 fresh rows right above the old row 106 (the "2010-19" block).
# Excel copies the formatting of the row above (row 104) onto the new rows,
# which is what gives the new row 105 its number formats and leaves rows
# 106:109 as blank-but-formatted placeholder rows (row 110 stays a true
# empty gap row, matching the existing blank-row rhythm used throughout
# this sheet).
$ws.Rows("105:109").Insert()

# New results row for the "2010-18" baseline re-run with wetland 10 active.
$ws.Range("A105").Value = "CW3M 1.1.0"
$ws.Range("B105").Value = "Demo_Baseline_2010-18 C600"
$ws.Range("C105").Value = "2010-18"
$ws.Range("D105").Value = 929.46866188888873
$ws.Range("E105").Value = 1890.2624918888889
$ws.Range("F105").Value = 1.0681051111111111
$ws.Range("G105").Value = 270.41205844444437
$ws.Range("H105").Value = 9.8445367777777779
$ws.Range("I105").Value = 7.3212358888888884
$ws.Range("J105").Value = 8.2027718888888881
$ws.Range("K105").Value = 668.60776777777789
$ws.Range("L105").Value = 80.524254777777799
$ws.Range("M105").Value = 1418.8558755555559
$ws.Range("N105").Value = 932.34357366666654
$ws.Range("O105").Value = 5824.0346137777778
$ws.Range("P105").Value = 27412.728515555555
$ws.Range("Q105").Value = 0.15715366666666666
$ws.Range("R105").Value = 0.000028666666666666671

# Match the author's final selection/scroll position on the sheet.
$ws.Activate()
$ws.Range("B106").Select()
